# Applies the Afrikaans translation edits described by the commit diff.
# Uses Find (search-only) to locate each English run, then assigns the
# Afrikaans replacement directly onto the found Range.Text. This sidesteps
# Find.Execute's built-in ReplaceWith path, which smart-quotes straight
# apostrophes (') into curly ones - the source diff uses plain ASCII quotes.

$d = $word.ActiveDocument
$missing = New-Object System.Collections.ArrayList

$rng = $d.Content
$found = $rng.Find.Execute("Appendix 16: SWIFT Interview Guide: Engagement  ", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found) {
    $rng.Text = "Bylaag 16: SWIFT Onderhoudgids: Betrokkenheid  "
} else {
    $missing.Add("Appendix 16: SWIFT Interview Guide: Engagement  ") | Out-Null
}

$rng = $d.Content
$found = $rng.Find.Execute("Briefing:", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found) {
    $rng.Text = "Inligtingsessie:"
} else {
    $missing.Add("Briefing:") | Out-Null
}

$rng = $d.Content
$found = $rng.Find.Execute("Hi there. Thank you for making the time for this phone call; it won’t take longer than 15 minutes. We noticed that you didn’t finish working through the ParentText programme. We would like to hear a little bit more about your experience, so we can hopefully improve the programme. ", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found) {
    $rng.Text = "Hi daar. Dankie dat jy tyd gemaak het vir hierdie oproep; dit sal nie langer as 15 minute neem nie. Ons het opgemerk dat jy nie die ParentText-program voltooi het nie. Ons sal graag meer wil hoor oor jou ervaring, sodat ons hopelik die program kan verbeter. "
} else {
    $missing.Add("Hi there. Thank you for making the time for this phone call; it won’t take longer than 15 minutes. We noticed that you didn’t finish working through the ParentText programme. We would like to hear a little bit more about your experience, so we can hopefully improve the programme. ") | Out-Null
}

$rng = $d.Content
$found = $rng.Find.Execute("We will be recording this phone call, so we can remember your answers, but your name and your answers will be kept confidential and will only be viewed by the research team. After the interview, the recording will be transcribed, or written down, and then it will be deleted. The written down information will be saved securely on a password protected computer. Do I have your permission to record the interview? ", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found) {
    $rng.Text = "Ons gaan hierdie telefoongesprek opneem, sodat ons jou antwoorde kan onthou, maar jou naam en antwoorde sal vertroulik gehou word en slegs deur die navorsingspan besigtig word. Na die onderhoud sal die opname transkribeer of neergeskryf word, en dan sal dit verwyder word. Die inligting wat neergeskryf word, sal veilig op ’n wagwoordbeskermde rekenaar gestoor word. Gee jy toestemming dat ons die onderhoud opneem? "
} else {
    $missing.Add("We will be recording this phone call, so we can remember your answers, but your name and your answers will be kept confidential and will only be viewed by the research team. After the interview, the recording will be transcribed, or written down, and then it will be deleted. The written down information will be saved securely on a password protected computer. Do I have your permission to record the interview? ") | Out-Null
}

$rng = $d.Content
$found = $rng.Find.Execute("There are no right or wrong answers. You can skip any questions you do not feel comfortable answering. You can also stop this conversation at any time if you wish. If you decide at a later stage that you would like your contribution to be removed from the study, you can contact the research team by email until the [date to be determined]. ", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found) {
    $rng.Text = "Daar is geen regte of verkeerde antwoorde nie. Jy kan enige vrae oorslaan waarmee jy nie gemaklik voel nie. Jy kan ook die gesprek op enige tyd stop as jy wil. As jy op 'n later stadium besluit dat jy jou bydrae uit die studie wil laat verwyder, kan jy die navorsingspan per e-pos kontak tot die [datum moet nog bepaal word]. "
} else {
    $missing.Add("There are no right or wrong answers. You can skip any questions you do not feel comfortable answering. You can also stop this conversation at any time if you wish. If you decide at a later stage that you would like your contribution to be removed from the study, you can contact the research team by email until the [date to be determined]. ") | Out-Null
}

$rng = $d.Content
$found = $rng.Find.Execute("Do you understand what I’ve just explained? Do you have any questions? Can we begin?", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found) {
    $rng.Text = "Verstaan jy wat ek sopas verduidelik het? Het jy enige vrae? Kan ons begin?"
} else {
    $missing.Add("Do you understand what I’ve just explained? Do you have any questions? Can we begin?") | Out-Null
}

$rng = $d.Content
$found = $rng.Find.Execute("What are the factors that you think led to you not finishing the programme?", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found) {
    $rng.Text = "Wat is die faktore wat jy dink daartoe gelei het dat jy nie die program voltooi het nie?"
} else {
    $missing.Add("What are the factors that you think led to you not finishing the programme?") | Out-Null
}

$rng = $d.Content
$found = $rng.Find.Execute("What was your experience of the content?", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found) {
    $rng.Text = "Wat was jou ervaring van die inhoud?"
} else {
    $missing.Add("What was your experience of the content?") | Out-Null
}

$rng = $d.Content
$found = $rng.Find.Execute("Probe: Explore relevance of content", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found) {
    $rng.Text = "Ondersoek: Verken die relevansie van die inhoud"
} else {
    $missing.Add("Probe: Explore relevance of content") | Out-Null
}

$rng = $d.Content
$found = $rng.Find.Execute("Probe: Explore interest in content", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found) {
    $rng.Text = "Ondersoek: Verken die belangstelling in die inhoud"
} else {
    $missing.Add("Probe: Explore interest in content") | Out-Null
}

$rng = $d.Content
$found = $rng.Find.Execute("Probe: What content could we have added to improve your experience? ", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found) {
    $rng.Text = "Ondersoek: Watter inhoud sou ons kon byvoeg om jou ervaring te verbeter? "
} else {
    $missing.Add("Probe: What content could we have added to improve your experience? ") | Out-Null
}

$rng = $d.Content
$found = $rng.Find.Execute("What was your experience of the time and data needed to complete the programme?", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found) {
    $rng.Text = "Wat was jou ervaring van die tyd en data wat nodig was om die program te voltooi?"
} else {
    $missing.Add("What was your experience of the time and data needed to complete the programme?") | Out-Null
}

$rng = $d.Content
$found = $rng.Find.Execute("Probe for data: Explore the use of Wi-Fi hotspots in the community and what they used to connect to ParentText; explore customisation of content delivery e.g. whether they used audio/visual only.", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found) {
    $rng.Text = "Ondersoek vir data: Verken die gebruik van Wi-Fi hotspots in die gemeenskap en wat hulle gebruik het om met ParentText te koppel; verken die aanpassing van inhoudlewering bv. of hulle net audio/visueel gebruik het."
} else {
    $missing.Add("Probe for data: Explore the use of Wi-Fi hotspots in the community and what they used to connect to ParentText; explore customisation of content delivery e.g. whether they used audio/visual only.") | Out-Null
}

$rng = $d.Content
$found = $rng.Find.Execute("Probe for time: Explore the amount and timing of messages being sent by the chatbot.", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found) {
    $rng.Text = "Ondersoek vir tyd: Verken die hoeveelheid en tydsberekening van die boodskappe wat deur die geselsbot gestuur is."
} else {
    $missing.Add("Probe for time: Explore the amount and timing of messages being sent by the chatbot.") | Out-Null
}

$rng = $d.Content
$found = $rng.Find.Execute("What was your experience of the home exercises/activities, and do you feel like they impacted whether you finished the programme?", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found) {
    $rng.Text = "Wat was jou ervaring van die tuisoefeninge/aktiwiteite, en voel jy dat hulle 'n impak gehad het of jy die program voltooi het?"
} else {
    $missing.Add("What was your experience of the home exercises/activities, and do you feel like they impacted whether you finished the programme?") | Out-Null
}

$rng = $d.Content
$found = $rng.Find.Execute("What do you think about the need for a programme like this to support parents?", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found) {
    $rng.Text = "Wat dink jy van die behoefte aan 'n program soos hierdie om ouers te ondersteun?"
} else {
    $missing.Add("What do you think about the need for a programme like this to support parents?") | Out-Null
}

$rng = $d.Content
$found = $rng.Find.Execute("What else can you recommend we do to improve the programme?", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found) {
    $rng.Text = "Wat anders kan jy aanbeveel wat ons kan doen om die program te verbeter?"
} else {
    $missing.Add("What else can you recommend we do to improve the programme?") | Out-Null
}

$rng = $d.Content
$found = $rng.Find.Execute("Is there anything that we haven’t spoken about that you’d like us to know? If after this conversation, there is anything that causes you to worry, remember you can still access the referral services in the programme by typing ‘help’. SADAG might be a good option. If you have any other questions about the study, you can send us a WhatsApp or email. ", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found) {
    $rng.Text = "Is daar iets wat ons nog nie oor gepraat het nie, wat jy graag wil hê ons moet weet? As daar iets is wat jou na hierdie gesprek bekommer, onthou jy kan steeds toegang kry tot die verwysingsdienste in die program deur ‘help’ te tik. SADAG mag dalk ‘n goeie opsie wees. As jy enige ander vrae oor die studie het, kan jy vir ons ‘n WhatsApp of e-pos stuur. "
} else {
    $missing.Add("Is there anything that we haven’t spoken about that you’d like us to know? If after this conversation, there is anything that causes you to worry, remember you can still access the referral services in the programme by typing ‘help’. SADAG might be a good option. If you have any other questions about the study, you can send us a WhatsApp or email. ") | Out-Null
}

if ($missing.Count -gt 0) {
    Write-Output "MISSING: $($missing.Count) segment(s) not found"
    foreach ($m in $missing) { Write-Output $m }
} else {
    Write-Output "All segments replaced successfully"
}
